$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style from an existing header cell (A1) onto the new
# header cell F1, then set its text.
$ws.Range("A1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$ws.Range("F1").Value = "Balanced_total"

# Balanced_total values for rows 2-20 (column F)
$values = @{
    2  = 11408
    3  = 548
    4  = 2440
    5  = 110
    6  = 356
    7  = 242
    8  = 770
    9  = 5506
    10 = 13336
    11 = 13976
    12 = 2708
    13 = 330
    14 = 1818
    15 = 1092
    16 = 6726
    17 = 784
    18 = 2372
    19 = 15398
    20 = 9346
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 6).Value = $values[$row]
}
